$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.933.37"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.089.00"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "2.398.64"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.772"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "2.083.16"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "37.936.23"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.21%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.14%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.31%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +6.94%  "
$ws.Range("D41").Value = "1.544.76"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.44%  "
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").Value = "2.284.96"
$ws.Range("E51").Value = "  +2.80%  "
